{"js": "// Neutralize racial-classification language to \"50M voters\" in three spots:\n//   1) Professional summary paragraph (plain text, no run split).\n//   2) \"Discovered systematic race coding errors...\" bullet \u2014 here \"50M\"\n//      must become its own bold/colored run (matching the \"23%\"/\"64%\" runs\n//      that already exist later in the same paragraph).\n//   3) \"Impact: Corrected demographic data...\" paragraph (plain text).\n\n// 1) Professional summary.\nconst summaryHits = context.document.body.search(\n  \"affecting all Black and Asian-American voters, developed geospatial ML\",\n  { matchCase: true }\n);\nsummaryHits.load(\"items\");\nawait context.sync();\nif (summaryHits.items.length > 0) {\n  summaryHits.items[0].insertText(\n    \"affecting 50M voters, developed geospatial ML\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) Bullet under \"Partner - Siege Analytics\": grab the containing paragraph\n//    first (so the follow-up \"50M\" lookup can be scoped to just this\n//    paragraph, rather than the whole document body), replace the plain\n//    phrase (keeping everything in one run), then re-find the freshly\n//    inserted \"50M\" token *within that paragraph* and flip its font to\n//    bold + the same dark slate color used by the other stat callouts in\n//    this paragraph (\"23%\", \"64%\", etc.).\nconst bulletHits = context.document.body.search(\n  \"affecting all Black and Asian-American voters, developed geospatial machine learning\",\n  { matchCase: true }\n);\nbulletHits.load(\"items\");\nawait context.sync();\nif (bulletHits.items.length > 0) {\n  const hit = bulletHits.items[0];\n  const bulletPara = hit.paragraphs.getFirst();\n  await context.sync();\n\n  hit.insertText(\n    \"affecting 50M voters, developed geospatial machine learning\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  const boldHits = bulletPara.search(\"50M\", { matchCase: true });\n  boldHits.load(\"items\");\n  await context.sync();\n  if (boldHits.items.length > 0) {\n    const target = boldHits.items[0];\n    target.font.bold = true;\n    target.font.color = \"#2C3E50\";\n    await context.sync();\n  }\n}\n\n// 3) \"Impact: Corrected demographic data...\" project-impact paragraph.\nconst impactHits = context.document.body.search(\n  \"affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%\",\n  { matchCase: true }\n);\nimpactHits.load(\"items\");\nawait context.sync();\nif (impactHits.items.length > 0) {\n  impactHits.items[0].insertText(\n    \"affecting 50M voters nationwide, improved electoral prediction accuracy by 22%\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Neutralize racial-classification language to \"50M voters\" in three spots:\n#   1) Professional summary paragraph (plain text, no run split).\n#   2) \"Discovered systematic race coding errors...\" bullet - here \"50M\"\n#      must become its own bold/colored run (matching the \"23%\"/\"64%\" runs\n#      that already exist later in the same paragraph).\n#   3) \"Impact: Corrected demographic data...\" paragraph (plain text).\n\n$d = $word.ActiveDocument\n\n# Dark-slate accent color used elsewhere in the resume for bolded stats\n# (hex 2C3E50 -> Word's BGR-packed wdColor long).\n$accentColor = 44 + (62 * 256) + (80 * 65536)\n\n# 1) Professional summary: plain substring swap, single run stays intact.\n$f1 = $d.Content.Find\n$f1.ClearFormatting()\n$f1.Replacement.ClearFormatting()\n$f1.Execute(\"affecting all Black and Asian-American voters, developed geospatial ML\", $true, $false, $false, $false, $false, $true, 1, $false, \"affecting 50M voters, developed geospatial ML\", 2)\n\n# 2) Bullet under \"Partner - Siege Analytics\": locate the containing\n#    paragraph first so the follow-up \"50M\" lookup can be scoped to just\n#    this paragraph instead of the whole document.\n$bulletPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Discovered systematic race coding errors affecting all Black and Asian-American voters*\") {\n        $bulletPara = $p\n        break\n    }\n}\n\n$f2 = $bulletPara.Range.Find\n$f2.ClearFormatting()\n$f2.Replacement.ClearFormatting()\n$f2.Execute(\"affecting all Black and Asian-American voters, developed geospatial machine learning\", $true, $false, $false, $false, $false, $true, 1, $false, \"affecting 50M voters, developed geospatial machine learning\", 2)\n\n# Re-find \"50M\" scoped to the same paragraph and promote it to its own\n# bold + colored run, matching the \"23%\"/\"64%\" callouts later on.\n$boldRng = $bulletPara.Range\n$f3 = $boldRng.Find\n$f3.ClearFormatting()\n$f3.Execute(\"50M\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n$boldRng.Font.Bold = 1\n$boldRng.Font.Color = $accentColor\n\n# 3) \"Impact: Corrected demographic data...\" project-impact paragraph.\n$f4 = $d.Content.Find\n$f4.ClearFormatting()\n$f4.Replacement.ClearFormatting()\n$f4.Execute(\"affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%\", $true, $false, $false, $false, $false, $true, 1, $false, \"affecting 50M voters nationwide, improved electoral prediction accuracy by 22%\", 2)\n"}
